# Daily attendance processing - reverse order of names in the
# "Recorded By" column (G) for every data row, except where the
# entry contains "admin@admin.com" (those rows are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2

    if ($v -eq $null) {
        continue
    }

    if ($v -like "*admin@admin.com*") {
        continue
    }

    $parts = $v -split ", "
    $n = $parts.Length

    if ($n -lt 2) {
        continue
    }

    $rev = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $rev += $parts[$i]
    }

    $newVal = $rev -join ", "

    if ($newVal -ne $v) {
        $cell.Value2 = $newVal
    }
}
